$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as per upstream GitHub Actions refresh.
# Price (column D) and Volume(1h) (column E) values are stored as TEXT in the source
# workbook (not numbers), including values that look numeric (e.g. "5.36"). When such a
# value is assigned directly to .Value, Excel auto-converts it to a Number, which would
# change the cell type. To preserve the original text semantics for those cells, we
# temporarily force a Text number format, assign the value, then restore the default
# "Normal" cell style so no visible formatting/style change is introduced.

$ws.Range("D2").Value = "68.983.66"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.744.97"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "2.743.54"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.365"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").Value = "3.243.22"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "68.938.78"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "2.711.59"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").Value = "2.878.66"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "602.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.79%  "
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("E49").Value = "  +6.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
